$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text runs collapse to plain text; all runs
#     shared identical formatting so this is visually equivalent) ---
$ws.Range("A8").Value = "Volume 30   Number  25"
$ws.Range("C9").Value = "Report Covering the Week  6/19/2023  Through  6/25/2023"

# --- Cells that flip between a numeric value and a text placeholder ("0" / "***.*") ---
# Pull the cell style from a stable same-category cell via .Copy() BEFORE any
# value edits happen (so the source cells still hold their original style/value),
# then overwrite with the final value.
$ws.Range("C14").Copy($ws.Range("D15"))   # numeric -> text placeholder "0"
$ws.Range("E14").Copy($ws.Range("E15"))   # numeric -> text placeholder "***.*"
$ws.Range("C15").Copy($ws.Range("C22"))   # numeric -> text placeholder "0"
$ws.Range("I14").Copy($ws.Range("D22"))   # text placeholder -> numeric
$ws.Range("K14").Copy($ws.Range("E22"))   # text placeholder -> numeric
$ws.Range("C15").Copy($ws.Range("C26"))   # numeric -> text placeholder "0"

$ws.Range("D15").Value = "0"
$ws.Range("E15").Value = "***.*"
$ws.Range("C22").Value = "0"
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = -100
$ws.Range("C26").Value = "0"

# --- Remaining same-category numeric value updates ---
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = -85.714285714285
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 0
$ws.Range("N15").Value = -21.212121212121
$ws.Range("C16").Value = 10
$ws.Range("E16").Value = -28.571428571428
$ws.Range("F16").Value = 47
$ws.Range("G16").Value = 49
$ws.Range("H16").Value = -4.081632653061
$ws.Range("I16").Value = 269
$ws.Range("J16").Value = 244
$ws.Range("K16").Value = 10.245901639344
$ws.Range("L16").Value = 62.048192771084
$ws.Range("M16").Value = 23.963133640553
$ws.Range("N16").Value = -70.044543429844
$ws.Range("C17").Value = 23
$ws.Range("D17").Value = 14
$ws.Range("E17").Value = 64.285714285714
$ws.Range("F17").Value = 78
$ws.Range("G17").Value = 74
$ws.Range("H17").Value = 5.405405405405
$ws.Range("I17").Value = 457
$ws.Range("J17").Value = 413
$ws.Range("K17").Value = 10.653753026634
$ws.Range("L17").Value = 46.006389776357
$ws.Range("M17").Value = 112.558139534884
$ws.Range("N17").Value = -10.03937007874
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 60
$ws.Range("F18").Value = 23
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = 9.523809523809
$ws.Range("I18").Value = 153
$ws.Range("J18").Value = 167
$ws.Range("K18").Value = -8.383233532934
$ws.Range("L18").Value = 84.33734939759
$ws.Range("M18").Value = 53
$ws.Range("N18").Value = -77.095808383233
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = 28.571428571428
$ws.Range("F19").Value = 58
$ws.Range("H19").Value = -4.918032786885
$ws.Range("I19").Value = 323
$ws.Range("J19").Value = 345
$ws.Range("K19").Value = -6.376811594202
$ws.Range("L19").Value = -5
$ws.Range("M19").Value = 71.808510638297
$ws.Range("N19").Value = -6.10465116279
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = 28.571428571428
$ws.Range("F20").Value = 27
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = -3.571428571428
$ws.Range("I20").Value = 170
$ws.Range("J20").Value = 161
$ws.Range("K20").Value = 5.590062111801
$ws.Range("L20").Value = 95.402298850574
$ws.Range("M20").Value = 209.090909090909
$ws.Range("N20").Value = -44.262295081967
$ws.Range("C21").Value = 68
$ws.Range("D21").Value = 54
$ws.Range("E21").Value = 25.925925925925
$ws.Range("F21").Value = 236
$ws.Range("G21").Value = 236
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 1403
$ws.Range("J21").Value = 1351
$ws.Range("K21").Value = 3.849000740192
$ws.Range("L21").Value = 39.048562933597
$ws.Range("M21").Value = 77.820025348542
$ws.Range("N21").Value = -49.731279111429
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 57
$ws.Range("K22").Value = -52.631578947368
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = -66.666666666666
$ws.Range("F23").Value = 38
$ws.Range("G23").Value = 45
$ws.Range("H23").Value = -15.555555555555
$ws.Range("I23").Value = 250
$ws.Range("J23").Value = 199
$ws.Range("K23").Value = 25.628140703517
$ws.Range("L23").Value = 92.307692307692
$ws.Range("M23").Value = 74.825174825174
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 43
$ws.Range("E24").Value = -25.581395348837
$ws.Range("F24").Value = 109
$ws.Range("G24").Value = 170
$ws.Range("H24").Value = -35.882352941176
$ws.Range("I24").Value = 745
$ws.Range("J24").Value = 785
$ws.Range("K24").Value = -5.095541401273
$ws.Range("L24").Value = 17.507886435331
$ws.Range("M24").Value = 20.745542949756
$ws.Range("C25").Value = 31
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = 93.75
$ws.Range("F25").Value = 94
$ws.Range("G25").Value = 97
$ws.Range("H25").Value = -3.092783505154
$ws.Range("I25").Value = 545
$ws.Range("J25").Value = 485
$ws.Range("K25").Value = 12.371134020618
$ws.Range("L25").Value = 25.287356321839
$ws.Range("M25").Value = 2.443609022556
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = -20
$ws.Range("J26").Value = 22
$ws.Range("K26").Value = 59.090909090909
$ws.Range("L26").Value = 0
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 10
$ws.Range("G27").Value = 11
$ws.Range("H27").Value = -9.090909090909
$ws.Range("I27").Value = 64
$ws.Range("J27").Value = 40
$ws.Range("K27").Value = 60
$ws.Range("L27").Value = 48.837209302325
$ws.Range("C28").Value = 1
$ws.Range("E28").Value = -50
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = -25
$ws.Range("I28").Value = 14
$ws.Range("J28").Value = 32
$ws.Range("K28").Value = -56.25
$ws.Range("L28").Value = -56.25
$ws.Range("M28").Value = -46.153846153846
$ws.Range("N28").Value = -84.946236559139
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 1
$ws.Range("I29").Value = 14
$ws.Range("J29").Value = 25
$ws.Range("K29").Value = -44
$ws.Range("L29").Value = -50
$ws.Range("M29").Value = -26.315789473684
$ws.Range("N29").Value = -83.132530120481
